$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns D and E (is_locked_lbl, is_enabled_lbl), shifting
# order_by / rem left into D / E (from F / G).
$ws.Range("D1:E1").EntireColumn.Delete()
